$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the "2.1 Objetivos Generales" body paragraph ("El objetivo de este
# proyecto...") via Find, so we are not dependent on hard-coded paragraph
# indices.
# ---------------------------------------------------------------------------
$search = $d.Content
$found = $search.Find.Execute("El objetivo de este proyecto es proporcionar", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the 'Objetivos Generales' body paragraph"
}

$bodyPara = $search.Paragraphs(1)

# ---------------------------------------------------------------------------
# Step 1: finish the sentence - add the rest of the paragraph's text as a
# second run, right before the paragraph mark.
# ---------------------------------------------------------------------------
$bodyRange = $bodyPara.Range
$insertPoint = $bodyRange.End - 1
$continuation = "lemas que poseen las empresas de servicio técnico, principalmente se encuentra en reducir los tiempos en el ciclo de reparación del articulo ingresado, balancear la carga de trabajo para los técnicos, optimizar las tareas criticas en los procesos de la organización y obtener información en tiempo real de estos. Con lo anterior se busca mejorar la calidad de servicio que se entrega al cliente y aumentar la capacidad para una mayor demanda de estos."

$bodyRange.InsertAfter($continuation)

# Force the newly-inserted text into its own run (instead of it silently
# being absorbed into the preceding run) by round-tripping a character
# property over just the new text.
$newTextRange = $d.Range($insertPoint, $insertPoint + $continuation.Length)
$originalBold = $newTextRange.Font.Bold
$newTextRange.Font.Bold = 1
$newTextRange.Font.Bold = $originalBold

# ---------------------------------------------------------------------------
# Step 2: this paragraph's mark carries a stray single-underline formatting,
# and is immediately followed by an empty paragraph. Deleting this
# paragraph's own trailing mark merges it forward into the following
# (plain) paragraph's mark - this both drops the stray underline formatting
# and collapses away the extra empty paragraph, exactly mirroring the
# source edit.
# ---------------------------------------------------------------------------
$bodyRange = $bodyPara.Range
$markRange = $d.Range($bodyRange.End - 1, $bodyRange.End)
$markRange.Delete()

# ---------------------------------------------------------------------------
# Step 3: insert a new empty paragraph right after the "2.2 Objetivo
# Especifico" heading, which now immediately follows the paragraph we just
# edited.
# ---------------------------------------------------------------------------
$headingPara = $bodyPara.Next()
$headingPara.Range.InsertParagraphAfter()
